$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds free-form text values (e.g. "65.218.24", "0.997",
# "0.0000260"); force text format so Excel does not reinterpret them as numbers
# and strip formatting such as trailing/leading zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.218.24"
$ws.Range("E2").Value = "  -5.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.240.86"
$ws.Range("E3").Value = "  -6.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.55"
$ws.Range("E5").Value = "  -5.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.15"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -5.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.241.17"
$ws.Range("E9").Value = "  -6.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  -12.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.571"
$ws.Range("E11").Value = "  -7.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.99"
$ws.Range("E12").Value = "  -9.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -7.49%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.49"
$ws.Range("E14").Value = "  -6.12%  "
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "623.82"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.743.95"
$ws.Range("E16").Value = "  -7.19%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.011.17"
$ws.Range("E17").Value = "  -5.67%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.73"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.228.12"
$ws.Range("E20").Value = "  -7.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("E21").Value = "  -8.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -5.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.66"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.24"
$ws.Range("E24").Value = "  +7.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.81"
$ws.Range("E25").Value = "  -8.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.84"
$ws.Range("E26").Value = "  -9.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.62"
$ws.Range("E27").Value = "  -8.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  -6.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.45"
$ws.Range("E29").Value = "  -8.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.52"
$ws.Range("E30").Value = "  -8.48%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  -7.77%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.69"
$ws.Range("E32").Value = "  -9.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.87"
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.102"
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "526.94"
$ws.Range("E35").Value = "  -9.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.719.36"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.58"
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0710"
$ws.Range("E39").Value = "  -8.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  -8.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("E41").Value = "  -4.19%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.32"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -6.80%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.24"
$ws.Range("E44").Value = "  -13.39%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.325"
$ws.Range("E46").Value = "  -12.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0405"
$ws.Range("E47").Value = "  -7.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.994"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").Value = "  -9.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.23"
$ws.Range("E51").Value = "  +0.23%  "
